$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3
$ws.Range("A3").Value = "2022-08-31 13:08:05"
$ws.Range("B3").Value = "MER_PREVENTION_INHAGOIA_12"
$ws.Range("C3").Value = "MER PREVENTION"

# Update row 4
$ws.Range("A4").Value = "2022-08-31 13:08:05"
$ws.Range("B4").Value = "MER_PREVENTION_INHAGOIA_12"
$ws.Range("C4").Value = "MER PREVENTION"

# Update row 5
$ws.Range("A5").Value = "2022-08-31 13:08:06"
$ws.Range("B5").Value = "MER_PREVENTION_INHAGOIA_12"
$ws.Range("C5").Value = "MER PREVENTION"
$ws.Range("D5").Value = "Buscar valores para cada indicador: DSD PREP"

# Update row 6
$ws.Range("A6").Value = "2022-08-31 13:08:11"
$ws.Range("B6").Value = "MER_PREVENTION_INHAGOIA_12"
$ws.Range("C6").Value = "MER PREVENTION"
$ws.Range("D6").Value = "Buscar valores para cada indicador: DSD TB PREV"

# Add new row 7
$ws.Range("A7").Value = "2022-08-31 13:08:15"
$ws.Range("B7").Value = "MER_PREVENTION_INHAGOIA_12"
$ws.Range("C7").Value = "MER PREVENTION"
$ws.Range("D7").Value = "Buscar valores para cada indicador: DSD GEND GBV"
$ws.Range("E7").Value = "ok"

# Add new row 8
$ws.Range("A8").Value = "2022-08-31 13:08:16"
$ws.Range("B8").Value = "MER_PREVENTION_INHAGOIA_12"
$ws.Range("C8").Value = "MER PREVENTION"
$ws.Range("D8").Value = "Buscar valores para cada indicador: DSD FPINT SITE"
$ws.Range("E8").Value = "ok"
